$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting so
# numeric-looking values (e.g. "542.48", "0.994") are stored as text,
# matching the source data which is inline text, not numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '60.861.75'
$ws.Range("E2").Value = '  +3.59%  '

$ws.Range("D3").Value = '2.561.89'
$ws.Range("E3").Value = '  +3.85%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '542.48'
$ws.Range("E5").Value = '  +2.06%  '

$ws.Range("D6").Value = '146.37'
$ws.Range("E6").Value = '  +1.90%  '

$ws.Range("D7").Value = '0.994'
$ws.Range("E7").Value = '  -0.31%  '

$ws.Range("D8").Value = '0.573'
$ws.Range("E8").Value = '  +0.76%  '

$ws.Range("D9").Value = '2.575.02'
$ws.Range("E9").Value = '  +3.41%  '

$ws.Range("E10").Value = '  +2.78%  '

$ws.Range("E11").Value = '  +1.40%  '

$ws.Range("D12").Value = '5.50'
$ws.Range("E12").Value = '  -1.91%  '

$ws.Range("D13").Value = '0.364'
$ws.Range("E13").Value = '  +4.28%  '

$ws.Range("D14").Value = '2.995.13'
$ws.Range("E14").Value = '  +3.13%  '

$ws.Range("D15").Value = '24.42'
$ws.Range("E15").Value = '  +2.41%  '

$ws.Range("D16").Value = '60.664.37'
$ws.Range("E16").Value = '  +3.42%  '

$ws.Range("D17").Value = '0.0000145'
$ws.Range("E17").Value = '  +5.49%  '

$ws.Range("D18").Value = '2.564.61'
$ws.Range("E18").Value = '  +3.23%  '

$ws.Range("D19").Value = '11.36'
$ws.Range("E19").Value = '  +0.82%  '

$ws.Range("E20").Value = '  +2.04%  '

$ws.Range("D21").Value = '329.01'
$ws.Range("E21").Value = '  +2.24%  '

$ws.Range("E22").Value = '  +4.73%  '

$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.17%  '

$ws.Range("D24").Value = '63.44'
$ws.Range("E24").Value = '  +4.64%  '

$ws.Range("D25").Value = '0.441'
$ws.Range("E25").Value = '  +1.34%  '

$ws.Range("E26").Value = '  +4.18%  '

$ws.Range("D27").Value = '0.992'
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("D28").Value = '8.08'
$ws.Range("E28").Value = '  +5.04%  '

$ws.Range("D29").Value = '7.15'
$ws.Range("E29").Value = '  +4.24%  '

$ws.Range("D30").Value = '0.0₃0812'
$ws.Range("E30").Value = '  +5.71%  '

$ws.Range("D31").Value = '1.82'
$ws.Range("E31").Value = '  +2.56%  '

$ws.Range("D32").Value = '1.21'
$ws.Range("E32").Value = '  -1.15%  '

$ws.Range("D33").Value = '164.54'
$ws.Range("E33").Value = '  +4.22%  '

$ws.Range("E34").Value = '  +5.92%  '

$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("D36").Value = '18.87'
$ws.Range("E36").Value = '  +2.05%  '

$ws.Range("D37").Value = '4.49'
$ws.Range("E37").Value = '  +3.28%  '

$ws.Range("E38").Value = '  +3.68%  '

$ws.Range("D39").Value = '5.72'
$ws.Range("E39").Value = '  -0.51%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '37.13'
$ws.Range("E40").Value = '  +1.55%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").Value = '304.75'
$ws.Range("E41").Value = '  +0.76%  '

$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '3.77'
$ws.Range("E42").Value = '  +2.05%  '

$ws.Range("B43").Value = 'SuiNetwork'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D43").Value = '0.845'
$ws.Range("E43").Value = '  +6.52%  '

$ws.Range("D44").Value = '0.611'
$ws.Range("E44").Value = '  +3.40%  '

$ws.Range("D45").Value = '0.996'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").Value = '10.87'
$ws.Range("E46").Value = '  +0.81%  '

$ws.Range("D47").Value = '126.63'
$ws.Range("E47").Value = '  +2.52%  '

$ws.Range("D48").Value = '19.18'
$ws.Range("E48").Value = '  +3.99%  '

$ws.Range("D49").Value = '0.0942'
$ws.Range("E49").Value = '  +2.11%  '

$ws.Range("D50").Value = '0.0526'
$ws.Range("E50").Value = '  +1.61%  '

$ws.Range("D51").Value = '0.0233'
$ws.Range("E51").Value = '  +2.60%  '
